$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 14 new rows above the existing SAMPLE block (rows 1-5 -> 15-19),
# making room for the new SPACE / PROJECT / EXPERIMENT blocks.
$ws.Rows("1:14").Insert()

$ws.Range("A1").Value = "SPACE"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 13

$ws.Range("B1").Value = ""
$ws.Range("B1").HorizontalAlignment = -4131

$ws.Range("A2").Value = "Code"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Font.Size = 13
$ws.Range("A2").HorizontalAlignment = -4131

$ws.Range("B2").Value = "Description"
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Size = 13

$ws.Range("A3").Value = "ELN_SETTINGS"
$ws.Range("A3").HorizontalAlignment = -4131

$ws.Range("B3").Value = "ELN Settings Updated"

$ws.Range("A5").Value = "PROJECT"
$ws.Range("A5").Font.Bold = $true
$ws.Range("A5").Font.Size = 13

$ws.Range("B5").Value = ""
$ws.Range("B5").HorizontalAlignment = -4131

$ws.Range("A6").Value = "Code"
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").Font.Size = 13
$ws.Range("A6").HorizontalAlignment = -4131

$ws.Range("B6").Value = "Description"
$ws.Range("B6").Font.Bold = $true
$ws.Range("B6").Font.Size = 13

$ws.Range("C6").Value = "Space"
$ws.Range("C6").Font.Bold = $true
$ws.Range("C6").Font.Size = 13

$ws.Range("A7").Value = "DEFAULT_PROJECT"
$ws.Range("A7").HorizontalAlignment = -4131

$ws.Range("C7").Value = "ELN_SETTINGS"
$ws.Range("C7").HorizontalAlignment = -4131

$ws.Range("A9").Value = "EXPERIMENT"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Size = 13

$ws.Range("B9").Value = ""
$ws.Range("B9").HorizontalAlignment = -4131

$ws.Range("A10").Value = "Experiment type"
$ws.Range("A10").Font.Bold = $true
$ws.Range("A10").Font.Size = 13

$ws.Range("B10").Value = ""
$ws.Range("B10").Font.Bold = $true
$ws.Range("B10").Font.Size = 13
$ws.Range("B10").Font.Color = 255

$ws.Range("C10").Value = ""
$ws.Range("C10").Font.Bold = $true
$ws.Range("C10").Font.Size = 13
$ws.Range("C10").Font.Color = 255

$ws.Range("A11").Value = "DEFAULT_EXPERIMENT"

$ws.Range("B11").Value = ""
$ws.Range("B11").HorizontalAlignment = -4131

$ws.Range("C11").Value = ""
$ws.Range("C11").HorizontalAlignment = -4131

$ws.Range("A12").Value = "Code"
$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").Font.Size = 13
$ws.Range("A12").HorizontalAlignment = -4131

$ws.Range("B12").Value = "Project"
$ws.Range("B12").Font.Bold = $true
$ws.Range("B12").Font.Size = 13

$ws.Range("C12").Value = "Name"
$ws.Range("C12").Font.Bold = $true
$ws.Range("C12").Font.Size = 13
$ws.Range("C12").HorizontalAlignment = -4131

$ws.Range("A13").Value = "DEFAULT_EXPERIMENT"
$ws.Range("A13").HorizontalAlignment = -4131

$ws.Range("B13").Value = "/ELN_SETTINGS/DEFAULT_PROJECT"
$ws.Range("B13").HorizontalAlignment = -4131

$ws.Range("B7").Value = "Default Project Updated"

$ws.Range("C13").Value = "Default Experiment Updated"

# Column widths (character units derived to best match the recorded pixel widths)
$ws.Columns("A").ColumnWidth = 21.5
$ws.Columns("B").ColumnWidth = 28.833333333333336
$ws.Columns("C").ColumnWidth = 24.333333333333336

# Row heights for rows whose content uses the larger 13pt bold font
$ws.Rows("1").RowHeight = 17
$ws.Rows("2").RowHeight = 17
$ws.Rows("5").RowHeight = 17
$ws.Rows("6").RowHeight = 17
$ws.Rows("9").RowHeight = 17
$ws.Rows("10").RowHeight = 17
$ws.Rows("12").RowHeight = 17
$ws.Rows("15").RowHeight = 17
$ws.Rows("16").RowHeight = 17
$ws.Rows("18").RowHeight = 17

# Restore the active selection
$ws.Range("C6").Select()
